# Plantilla Lista de Tareas de la Entrega 2 - Entrega 2 edits
# Commit: "Diagramas de secuencia y registros de tiempo"
#   - Se agregaron los diagramas de secuencia de los CU 14,15,16,21
#   - Se registraron los tiempos que se invirtieron en realizarlos.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Mark the 4 "Diagrama de robustez y Diagrama de secuencia" tasks
#     (CU-14 / CU-15 / CU-16 / CU-21, rows 59/61/63/65) as done, and log the
#     hours consumed on the corresponding days (columns T="Dia 11",
#     Z="Dia 13", AC="Dia 14").

# CU-14 (row 59)
$ws.Range("F59").Value = "Hecho"
$ws.Range("T59").Value = 0.58
$ws.Range("Z59").Value = 1
$ws.Range("AA59").Value = 0.5

# CU-15 (row 61)
$ws.Range("F61").Value = "Hecho"
$ws.Range("T61").Value = 0.58
$ws.Range("Z61").Value = 0.58

# CU-16 (row 63)
$ws.Range("F63").Value = "Hecho"
$ws.Range("AC63").Value = 0.58

# CU-21 (row 65)
$ws.Range("F65").Value = "Hecho"
$ws.Range("AC65").Value = 0.58

# --- Narrow column C (description) now that the sequence diagrams are done
#     and the column is no longer needed as wide; this makes row 66's
#     wrapped text take up more lines.
$ws.Columns("C").ColumnWidth = 58.14
$ws.Rows(66).RowHeight = 38.25

# --- Update the window view state to reflect where the user ended up.
$ws.Range("D61").Select()
$excel.ActiveWindow.Zoom = 100
